$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("Premium") model list - wrap model #1 (KANNON PAINTED) in a hyperlink
$ws.Range("B16").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/kannon-painted/'' target=''_blank''>/KANNON PAINTED</a> 2. TRINTON PRIME 3. AVANCER E-SENSE 4. TRINTON ELECTROPLATED 5. TRINTON PAINTED 6. STUBBY 7. RIVERA 600 MM 8. AGGNUZ PRIMUS 9. TWIRL 10. RADIANCE 11. AVANCER PRIME-ANTI DUST 12. AVANCER 13. AMOUR 14. IMPERIAL 15. LERONE 16. SENORITA 17. SENORITA 1200 18. DEVINE 19. AURA PRIME - ANTI DUST'

# Row 17 ("Metallic / Decorative") model list - wrap model #1 (AIR 360 DECO) in a hyperlink
$ws.Range("B17").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/air-360-deco/" target=''_blank''>/AIR 360 DECO</a> 2. CALIBER (HI SPEED) 3. AURA 4. RIVIERA STYLE 5. FLYLEAF 6. DECAIR 7. CALIBER (HI SPEED) METALL 8. AURA METALLIC 9. HIGH SPEED DECORA METALLIC 10. HIGH SPEED DECORA 11. DIAMOND 12. FLAVIA 13. MOUNT BREESE DECO 14. HIGH SPEED METALLIC'

# Restore view: active cell back to B17, top-left reset to default
$ws.Range("B17").Select()
